$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the authoritative diff.
# Numeric-looking strings get NumberFormat "@" first so Excel
# stores them as text (matching the original inlineStr cells)
# instead of silently parsing them into floating point numbers.

$ws.Cells.Item(2, 4).Value = '43.786.68'
$ws.Cells.Item(2, 5).Value = '  -0.82%  '
$ws.Cells.Item(3, 4).Value = '2.235.33'
$ws.Cells.Item(3, 5).Value = '  -1.80%  '
$ws.Cells.Item(4, 5).Value = '  +0.16%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '313.71'
$ws.Cells.Item(5, 5).Value = '  -1.86%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '98.28'
$ws.Cells.Item(6, 5).Value = '  -5.20%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.571'
$ws.Cells.Item(7, 5).Value = '  -3.17%  '
$ws.Cells.Item(8, 5).Value = '  +0.28%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.532'
$ws.Cells.Item(9, 5).Value = '  -7.09%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '35.67'
$ws.Cells.Item(10, 5).Value = '  -8.24%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0821'
$ws.Cells.Item(11, 5).Value = '  -2.33%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '7.35'
$ws.Cells.Item(12, 5).Value = '  -6.93%  '
$ws.Cells.Item(13, 5).Value = '  -2.90%  '
$ws.Cells.Item(14, 4).Value = '2.575.28'
$ws.Cells.Item(14, 5).Value = '  -1.86%  '
$ws.Cells.Item(15, 2).Value = 'Polygon'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.839'
$ws.Cells.Item(15, 5).Value = '  -4.51%  '
$ws.Cells.Item(16, 2).Value = 'WrappedEther'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(16, 4).Value = '2.236.94'
$ws.Cells.Item(16, 5).Value = '  -1.80%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '13.99'
$ws.Cells.Item(17, 5).Value = '  -3.83%  '
$ws.Cells.Item(18, 4).Value = '43.640.12'
$ws.Cells.Item(18, 5).Value = '  -1.18%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '13.10'
$ws.Cells.Item(19, 5).Value = '  -7.42%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0965'
$ws.Cells.Item(20, 5).Value = '  -3.42%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '6.30'
$ws.Cells.Item(21, 5).Value = '  -4.92%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '65.54'
$ws.Cells.Item(22, 5).Value = '  -1.04%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '235.77'
$ws.Cells.Item(23, 5).Value = '  -0.80%  '
$ws.Cells.Item(24, 5).Value = '  -7.21%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.02'
$ws.Cells.Item(25, 5).Value = '  -8.14%  '
$ws.Cells.Item(26, 5).Value = '  +0.15%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.09'
$ws.Cells.Item(27, 5).Value = '  -1.70%  '
$ws.Cells.Item(28, 5).Value = '  -3.16%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '36.51'
$ws.Cells.Item(29, 5).Value = '  -6.54%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '5.96'
$ws.Cells.Item(30, 5).Value = '  -8.96%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '156.79'
$ws.Cells.Item(31, 5).Value = '  -3.50%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '19.91'
$ws.Cells.Item(32, 5).Value = '  -2.94%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.0833'
$ws.Cells.Item(33, 5).Value = '  -5.78%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '2.65'
$ws.Cells.Item(34, 5).Value = '  -2.28%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '3.28'
$ws.Cells.Item(35, 5).Value = '  +0.56%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.90'
$ws.Cells.Item(36, 5).Value = '  -8.49%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.109'
$ws.Cells.Item(37, 5).Value = '  -2.39%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.117'
$ws.Cells.Item(38, 5).Value = '  -3.68%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '15.64'
$ws.Cells.Item(39, 5).Value = '  +1.56%  '
$ws.Cells.Item(40, 2).Value = 'NEARProtocol'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.54'
$ws.Cells.Item(40, 5).Value = '  -9.37%  '
$ws.Cells.Item(41, 2).Value = 'RenderToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '4.02'
$ws.Cells.Item(41, 5).Value = '  -10.92%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.0307'
$ws.Cells.Item(42, 5).Value = '  -5.94%  '
$ws.Cells.Item(43, 5).Value = '  +0.12%  '
$ws.Cells.Item(44, 4).Value = '1.704.09'
$ws.Cells.Item(44, 5).Value = '  -4.16%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '82.80'
$ws.Cells.Item(45, 5).Value = '  -2.91%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.194'
$ws.Cells.Item(46, 5).Value = '  -6.84%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '5.12'
$ws.Cells.Item(47, 5).Value = '  -4.53%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '101.43'
$ws.Cells.Item(48, 5).Value = '  -3.03%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '71.29'
$ws.Cells.Item(49, 5).Value = '  -4.90%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.62'
$ws.Cells.Item(50, 5).Value = '  +0.23%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '56.21'
$ws.Cells.Item(51, 5).Value = '  -6.02%  '
